# Edit: split the run containing "Os valores devem ser exibidos utilizando o
# formato R$ xxx.xx, exemplo:" into three runs, wrapping the "xxx.xx" token
# with <w:proofErr w:type="spellStart"/> / <w:proofErr w:type="spellEnd"/>
# markers, matching what Word's proofing pass produces when it flags
# "xxx.xx" as a possible spelling error.

$d = $word.ActiveDocument

$needle = "xxx.xx"

foreach ($p in $d.Paragraphs) {
    $paraText = $p.Range.Text
    if ($paraText -like "*$needle*") {

        # Pull the paragraph's own OOXML so we can reuse its real paragraph
        # mark attributes (paraId, rsids, etc.) and run/rPr formatting
        # instead of inventing new ones.
        $full = $p.Range.WordOpenXML

        $pOpen = "<w:p>"
        if ($full -match '(<w:p [^>]*>|<w:p>)') { $pOpen = $matches[1] }

        $pPr = ""
        if ($full -match '(<w:pPr>.*?</w:pPr>)') { $pPr = $matches[1] }

        $rOpenWithRsid = "<w:r>"
        if ($full -match '(<w:r [^>]*>|<w:r>)') { $rOpenWithRsid = $matches[1] }

        $rPr = ""
        if ($full -match '<w:r[^>]*>(<w:rPr>.*?</w:rPr>)') { $rPr = $matches[1] }

        # Text of the (single) run, split around the "xxx.xx" token.
        $runText = $p.Range.Text.TrimEnd([char]13, [char]7)
        $idx = $runText.IndexOf($needle)
        $before = $runText.Substring(0, $idx)
        $after = $runText.Substring($idx + $needle.Length)

        $xmlns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

        function Make-TextElement([string]$text) {
            if ($text.Length -gt 0 -and ($text[0] -eq ' ' -or $text[-1] -eq ' ')) {
                return "<w:t xml:space=`"preserve`">" + $text + "</w:t>"
            }
            return "<w:t>" + $text + "</w:t>"
        }

        $newXml = $pOpen + $pPr +
                  $rOpenWithRsid + $rPr + (Make-TextElement $before) + "</w:r>" +
                  "<w:proofErr w:type=`"spellStart`"/>" +
                  "<w:r $xmlns>" + $rPr + (Make-TextElement $needle) + "</w:r>" +
                  "<w:proofErr w:type=`"spellEnd`"/>" +
                  "<w:r $xmlns>" + $rPr + (Make-TextElement $after) + "</w:r>" +
                  "</w:p>"

        $p.Range.InsertXML($newXml)
    }
}
